# Trade #18 closed at 2026-02-17 13:17:56 - unknown UNKNOWN +0.000%
#
# Updates the "Summary" and "Strategy Status" aggregate figures to reflect
# the newly-closed trade, and appends the trade's detail row to both the
# "All Trades" and "MarketMaking" logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet - refresh aggregate stats
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.27    # Current Capital
$summary.Range("B4").Value = -0.73      # Total P&L $
$summary.Range("B5").Value = -0.81      # Total P&L %
$summary.Range("B6").Value = 18         # Total Trades
$summary.Range("B7").Value = 6          # Winning Trades
$summary.Range("B9").Value = 33.33      # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet - refresh the MarketMaking row
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.27       # Capital
$status.Range("D4").Value = 18          # Trades
$status.Range("E4").Value = -0.73       # P&L $
$status.Range("F4").Value = -0.73       # P&L %
$status.Range("G4").Value = 33.33       # Win Rate %

# ---------------------------------------------------------------
# Append the new trade record (row 19) to both trades-log worksheets.
# ---------------------------------------------------------------
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(19, 1).Value = 18                 # Trade #
    $ws.Cells.Item(19, 2).Value = "'2026-02-17"       # Date (force text)
    $ws.Cells.Item(19, 3).Value = "'13:17:49"         # Time (force text)
    $ws.Cells.Item(19, 4).Value = "MarketMaking"      # Strategy
    $ws.Cells.Item(19, 5).Value = "DOWN"              # Side
    $ws.Cells.Item(19, 6).Value = 0.81                # Entry Price
    $ws.Cells.Item(19, 7).Value = 0.83                # Exit Price
    $ws.Cells.Item(19, 8).Value = "CLOSED"            # Status
    $ws.Cells.Item(19, 9).Value = 2.4691              # P&L %
    $ws.Cells.Item(19, 10).Value = 0.02               # P&L $
    $ws.Cells.Item(19, 11).Value = 99.27              # Capital After
    $ws.Cells.Item(19, 12).Value = 0                  # Entry Slippage (bps)
    $ws.Cells.Item(19, 13).Value = 0                  # Exit Slippage (bps)
    $ws.Cells.Item(19, 14).Value = 0.6                # Confidence
    $ws.Cells.Item(19, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item(19, 16).Value = "early_exit"       # Exit Reason
    $ws.Cells.Item(19, 17).Value = 0.14               # Duration (min)
}
